$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '46.886.49'
$ws.Range('E2').Value = '  +4.78%  '
$ws.Range('D3').Value = '2.476.69'
$ws.Range('E3').Value = '  +2.25%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '323.87'
$ws.Range('E5').Value = '  +2.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '105.07'
$ws.Range('E6').Value = '  +3.52%  '
$ws.Range('E7').Value = '  +1.72%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.541'
$ws.Range('E9').Value = '  +1.98%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.18'
$ws.Range('E10').Value = '  +2.26%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0816'
$ws.Range('E11').Value = '  +2.02%  '
$ws.Range('E12').Value = '  +0.71%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.25'
$ws.Range('E13').Value = '  -2.64%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.13'
$ws.Range('E14').Value = '  +3.04%  '
$ws.Range('D15').Value = '2.863.19'
$ws.Range('E15').Value = '  +2.30%  '
$ws.Range('D16').Value = '2.453.70'
$ws.Range('E16').Value = '  +1.50%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.845'
$ws.Range('E17').Value = '  +1.56%  '
$ws.Range('D18').Value = '46.744.18'
$ws.Range('E18').Value = '  +4.97%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.62'
$ws.Range('E19').Value = '  +2.28%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.48'
$ws.Range('E20').Value = '  +1.69%  '
$ws.Range('D21').Value = '0.0₃0937'
$ws.Range('E21').Value = '  +2.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '70.52'
$ws.Range('E22').Value = '  +2.56%  '
$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '249.97'
$ws.Range('E23').Value = '  +3.07%  '
$ws.Range('B24').Value = 'ImmutableX'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.39'
$ws.Range('E24').Value = '  +4.80%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.56'
$ws.Range('E25').Value = '  +2.59%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.21'
$ws.Range('E26').Value = '  +4.00%  '
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('E28').Value = '  -1.92%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.85'
$ws.Range('E29').Value = '  +3.90%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '35.15'
$ws.Range('E30').Value = '  +4.45%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.136'
$ws.Range('E31').Value = '  +7.26%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '49.59'
$ws.Range('E32').Value = '  +2.20%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.67'
$ws.Range('E33').Value = '  +0.80%  '
$ws.Range('E34').Value = '  +2.91%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0768'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.62'
$ws.Range('E37').Value = '  +3.22%  '
$ws.Range('E38').Value = '  +1.53%  '
$ws.Range('E39').Value = '  +3.61%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '123.21'
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('E41').Value = '  +1.68%  '
$ws.Range('E42').Value = '  +1.10%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '20.98'
$ws.Range('E43').Value = '  -0.68%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0294'
$ws.Range('E44').Value = '  +0.97%  '
$ws.Range('D45').Value = '1.976.68'
$ws.Range('E45').Value = '  +1.83%  '
$ws.Range('E46').Value = '  +1.19%  '
$ws.Range('E47').Value = '  -0.62%  '
$ws.Range('E48').Value = '  +3.68%  '
$ws.Range('E49').Value = '  -2.57%  '
$ws.Range('E50').Value = '  +16.77%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '79.53'
$ws.Range('E51').Value = '  +4.97%  '
